$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cell value changes (row -> parameter)
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 7.5
$ws.Range("F2").Value = 2.8

$ws.Range("F3").Value = 86
$ws.Range("G3").Value = 96

$ws.Range("F4").Value = 86
$ws.Range("G4").Value = 96

$ws.Range("D5").Value = 56
$ws.Range("F5").Value = 66
$ws.Range("G5").Value = 92

$ws.Range("D6").Value = 50

$ws.Range("G8").Value = 95

$ws.Range("G9").Value = 2

$ws.Range("C10").Value = 0.13600000000000001
$ws.Range("D10").Value = 0.13600000000000001

$ws.Range("G11").Value = 0.00020000000000000001

$ws.Range("C16").Value = 1.8
$ws.Range("D16").Value = 1.22
$ws.Range("F16").Value = 0.87
$ws.Range("G16").Value = 0.76

# Update selection on sheet to G12
$ws.Range("G12").Select()
